# feat: add 2022-Q4 data
#
# Target end state:
#   Sheet 1: "总计"     (unchanged name/position) - gains a 2022-Q4 summary row
#   Sheet 2: "2022-Q4"  (was "2022-Q3") - now holds the new Q4 fund-holdings table
#   Sheet 3: "2022-Q3"  (new tab) - holds the fund-holdings table that used to live
#                         on the "2022-Q3" tab, moved here verbatim

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q3" tab (values + formats) so the data
#    currently on it survives, placed on a new tab right after it.
# ---------------------------------------------------------------------------
$q3.Copy($null, $q3)
$q3Moved = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 2) Rename tabs: the original tab becomes "2022-Q4" (overwritten with the
#    new quarter's data below); the duplicate becomes the permanent "2022-Q3".
# ---------------------------------------------------------------------------
$q4 = $q3
$q4.Name = "2022-Q4"
$q3Moved.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 3) Wipe the old contents of the (now-renamed) "2022-Q4" sheet and write the
#    brand-new fund-holdings table for 2022-Q4 in its place.
# ---------------------------------------------------------------------------
$q4.Cells.Clear()

# Mark the text-bearing ranges as Text ("@") BEFORE assigning values so
# Excel doesn't reinterpret numeric-looking strings (e.g. "009664", "7.30")
# as numbers and strip leading zeros / trailing zeros.
$q4.Range("B1:H1").NumberFormat = "@"
$q4.Range("B2:G7").NumberFormat = "@"

$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "009664"
$q4.Cells.Item(2,3).Value = "汇添富医疗积极成长一年持有期混合A"
$q4.Cells.Item(2,4).Value = "29.15"
$q4.Cells.Item(2,5).Value = "77.39"
$q4.Cells.Item(2,6).Value = "2.37"
$q4.Cells.Item(2,7).Value = "0.6909"
$q4.Cells.Item(2,8).Value = 10

$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "519171"
$q4.Cells.Item(3,3).Value = "浦银安盛医疗健康灵活配置混合A"
$q4.Cells.Item(3,4).Value = "7.30"
$q4.Cells.Item(3,5).Value = "92.50"
$q4.Cells.Item(3,6).Value = "2.95"
$q4.Cells.Item(3,7).Value = "0.2154"
$q4.Cells.Item(3,8).Value = 9

$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = "009665"
$q4.Cells.Item(4,3).Value = "汇添富医疗积极成长一年持有期混合C"
$q4.Cells.Item(4,4).Value = "3.62"
$q4.Cells.Item(4,5).Value = "77.39"
$q4.Cells.Item(4,6).Value = "2.37"
$q4.Cells.Item(4,7).Value = "0.0858"
$q4.Cells.Item(4,8).Value = 10

$q4.Cells.Item(5,1).Value = 3
$q4.Cells.Item(5,2).Value = "013183"
$q4.Cells.Item(5,3).Value = "浦银安盛医疗健康灵活配置混合C"
$q4.Cells.Item(5,4).Value = "1.08"
$q4.Cells.Item(5,5).Value = "92.50"
$q4.Cells.Item(5,6).Value = "2.95"
$q4.Cells.Item(5,7).Value = "0.0319"
$q4.Cells.Item(5,8).Value = 9

$q4.Cells.Item(6,1).Value = 4
$q4.Cells.Item(6,2).Value = "014547"
$q4.Cells.Item(6,3).Value = "财通医药鑫选6个月持有期混合A"
$q4.Cells.Item(6,4).Value = "0.44"
$q4.Cells.Item(6,5).Value = "89.44"
$q4.Cells.Item(6,6).Value = "4.23"
$q4.Cells.Item(6,7).Value = "0.0186"
$q4.Cells.Item(6,8).Value = 9

$q4.Cells.Item(7,1).Value = 5
$q4.Cells.Item(7,2).Value = "014548"
$q4.Cells.Item(7,3).Value = "财通医药鑫选6个月持有期混合C"
$q4.Cells.Item(7,4).Value = "0.28"
$q4.Cells.Item(7,5).Value = "89.44"
$q4.Cells.Item(7,6).Value = "4.23"
$q4.Cells.Item(7,7).Value = "0.0118"
$q4.Cells.Item(7,8).Value = 9

# Re-apply the proper (borderless / header) styling used elsewhere in this
# workbook: header row + index column use the same style as the "总计"
# sheet's header row; the plain text body cells get the default style back
# (PasteSpecial formats overwrites the "@" override applied above).
$summary.Cells.Item(1,2).Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$summary.Cells.Item(2,1).Copy()
$q4.Range("A2:A7").PasteSpecial(-4122)
$summary.Cells.Item(200,200).Copy()
$q4.Range("B2:G7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Update the "总计" summary sheet: push the existing 2022-Q3 total down to
#    row 3 (re-numbering its index column), and put the new 2022-Q4 totals
#    into row 2.
# ---------------------------------------------------------------------------
$summary.Cells.Item(2,1).Copy()
$summary.Cells.Item(3,1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q3"
$summary.Cells.Item(3,3).Value = 2
$summary.Cells.Item(3,4).Value = 0.05

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 6
$summary.Cells.Item(2,4).Value = 1.05
